$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.136116072538749
$ws.Range("C2").Value = 0.136116072538749
$ws.Range("D2").Value = 0.962610311559667
$ws.Range("E2").Value = 0.00528136954727986
$ws.Range("F2").Value = 0.448

$ws.Range("B3").Value = 0.428008694200286
$ws.Range("C3").Value = 0.428008694200286
$ws.Range("D3").Value = 3.02686945626568
$ws.Range("E3").Value = 0.016606944656568
$ws.Range("F3").Value = 0.0066

$ws.Range("B4").Value = 4.84432176042587
$ws.Range("C4").Value = 1.61477392014196
$ws.Range("D4").Value = 11.4196508713089
$ws.Range("E4").Value = 0.187962030828174

$ws.Range("B5").Value = 0.0818481532698945
$ws.Range("C5").Value = 0.0818481532698945
$ws.Range("D5").Value = 0.578828604515365
$ws.Range("E5").Value = 0.00317574799300543
$ws.Range("F5").Value = 0.7611

$ws.Range("B6").Value = 0.209823310614228
$ws.Range("C6").Value = 0.0699411035380759
$ws.Range("D6").Value = 0.494622172179173
$ws.Range("E6").Value = 0.00814124608739322
$ws.Range("F6").Value = 0.9675

$ws.Range("B7").Value = 1.2168096680377
$ws.Range("C7").Value = 0.405603222679233
$ws.Range("D7").Value = 2.86841838197845
$ws.Range("E7").Value = 0.0472128045259354
$ws.Range("F7").Value = 0.0003

$ws.Range("B8").Value = 0.190738581332644
$ws.Range("C8").Value = 0.0635795271108813
$ws.Range("D8").Value = 0.449633223024406
$ws.Range("E8").Value = 0.00740074934688417
$ws.Range("F8").Value = 0.9793

$ws.Range("B9").Value = 18.6652078825163
$ws.Range("C9").Value = 0.141403090019063
$ws.Range("E9").Value = 0.72421910701476

$ws.Range("B10").Value = 25.7728741229356
